# Removed the need of LocationId
#
# The first data column (A) held the "LocationId" (ns2:Sender / "ABC") value,
# which is no longer needed. Every other column (ns4:Base, ns4:Base3,
# ns3:ValidFrom, ns3:ValidTo, ns3:Quantity and their row values) shifts one
# column to the left (B->A, C->B, D->C, E->D, F->E) and the former last
# column (F) becomes empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 1
$lastRow  = 4
$firstCol = 2   # B
$lastCol  = 6   # F

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $src = $ws.Cells.Item($r, $c)
        $dst = $ws.Cells.Item($r, $c - 1)
        if ($src.Text -eq "") {
            $dst.Value = $null
        } else {
            $dst.Value = $src.Value2
        }
    }
    # the old last column (F) is now unused, clear it out
    $ws.Cells.Item($r, $lastCol).Clear() | Out-Null
}

# Update the selection / active cell as recorded in the saved workbook
$ws.Range("D12").Select() | Out-Null
